$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the keyword / appID pairs that changed.
$ws.Range("A8").Value = "optimistic quotes"

$ws.Range("A9").Value = "earn passive income"
$ws.Range("B9").Value = "passive.income.nadi.myfirstdrawermenuproject2"

$ws.Range("A10").Value = "duty free products"
$ws.Range("B10").Value = "duty.pare.myapp"

$ws.Range("A11").Value = "powerful quotes"
$ws.Range("B11").Value = "com.sugar.powerfulquotes"

$ws.Range("A12").Value = "duty free"
$ws.Range("B12").Value = "duty.pare.myapp"

$ws.Range("A13").Value = "motivation quotes"
$ws.Range("B13").Value = "com.sugar.powerfulquotes"

$ws.Range("A14").Value = "passive income"
$ws.Range("B14").Value = "passive.income.nadi.myfirstdrawermenuproject2"

$ws.Range("A15").Value = "duty free"
$ws.Range("B15").Value = "duty.pare.myapp"

$ws.Range("A16").Value = "motivation quotes"
$ws.Range("B16").Value = "com.sugar.powerfulquotes"

$ws.Range("A17").Value = "motivation quotes"
$ws.Range("B17").Value = "com.sugar.powerfulquotes"

$ws.Range("A18").Value = "optimistic quotes"

$ws.Range("A19").Value = "passive income"
$ws.Range("B19").Value = "passive.income.nadi.myfirstdrawermenuproject2"

$ws.Range("A20").Value = "duty free products"
$ws.Range("B20").Value = "duty.pare.myapp"

# Remove the trailing rows that no longer exist after the edit.
$ws.Rows("21:23").Delete()

# Match the author's final selection / scroll position.
$ws.Range("A9").Select() | Out-Null
